# Update the "修改时间" (last-modified timestamp) column on every portfolio
# sheet from 202509211902 to 202509212022.
#
# The timestamp is stored as TEXT (not a number) in the workbook. Excel's
# Range.Value setter auto-detects purely numeric strings and stores them as
# numbers, so we briefly force the cells to the Text number format before
# writing the value, then restore the cell style back to Normal so no extra
# formatting is left behind on the cells.

$wb = $excel.ActiveWorkbook
$newTimestamp = "202509212022"

# Sheet 1 "大智投资组合": timestamp lives in column E, data rows 2-9
$ws1 = $wb.Worksheets.Item(1)
$rng1 = $ws1.Range("E2:E9")
$rng1.NumberFormat = "@"
for ($row = 2; $row -le 9; $row++) {
    $ws1.Cells.Item($row, 5).Value = $newTimestamp
}
$rng1.Style = "Normal"

# Sheet 2 "大成投资组合": timestamp lives in column E, data rows 2-11
$ws2 = $wb.Worksheets.Item(2)
$rng2 = $ws2.Range("E2:E11")
$rng2.NumberFormat = "@"
for ($row = 2; $row -le 11; $row++) {
    $ws2.Cells.Item($row, 5).Value = $newTimestamp
}
$rng2.Style = "Normal"

# Sheet 3 "我的投资组合": timestamp lives in column G, data rows 2-13
$ws3 = $wb.Worksheets.Item(3)
$rng3 = $ws3.Range("G2:G13")
$rng3.NumberFormat = "@"
for ($row = 2; $row -le 13; $row++) {
    $ws3.Cells.Item($row, 7).Value = $newTimestamp
}
$rng3.Style = "Normal"
